# Apply the weekly work-report update:
#  - refresh "Report Generated On" timestamp
#  - update summary totals (Total Billed Amount, Total Line Items, Billing Period)
#  - insert a new billed line item (SVC-10-TP-AAA-RS) into the "Point 35 / Thursday" section
#  - update the grand TOTAL row accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header / summary fields
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:47 AM"

$ws.Range("C8").Value = 10974.66
$ws.Range("C9").Value = 81
$ws.Range("C10").Value = "06/30/2025 to 07/06/25"

# ---------------------------------------------------------------------------
# 2. Insert a new line item row above row 108 (pushes the existing
#    SVD-3-CV-C / SVD-SG2 / XCO-27-100-8-C rows and the TOTAL row down by one)
# ---------------------------------------------------------------------------
$ws.Rows("108:108").Insert()

# The detail rows use alternating row-banding (shaded / plain) keyed off the
# row number, not the data that happens to live in the row. Rows 103 (plain)
# and 104 (shaded) sit above the insertion point, so they are untouched by
# the shift and make reliable format sources. Re-stripe rows 108-111 so the
# banding keeps alternating correctly now that a row has been inserted
# (108 & 110 -> shaded like row 104, 109 & 111 -> plain like row 103).
$ws.Range("A104:H104").Copy()
$ws.Range("A108:H108").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A110:H110").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A103:H103").Copy()
$ws.Range("A109:H109").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A111:H111").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new line item
$ws.Range("A108").Value = "Point 35"
$ws.Range("B108").Value = "SVC-10-TP-AAA-RS"
$ws.Range("C108").Value = "Trans"
$ws.Range("D108").Value = "SVC,1/0,Trip,All Alum,Res"
$ws.Range("E108").Value = "EA"
$ws.Range("F108").Value = 1
$ws.Range("G108").Value = ""
$ws.Range("H108").Value = 195.83

# ---------------------------------------------------------------------------
# 3. Update the grand TOTAL row (now shifted from row 111 to row 112)
# ---------------------------------------------------------------------------
$ws.Range("H112").Value = 1881.26
